$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "Cookies"
$ws.Range("A1").Value = "Hello"
